$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shop candidate pointer: 1081 -> 822 ---
$ws.Range("G4").Value = 822

# --- Row 6 (shop item 1): title/notes now spans two lines ---
$ws.Range("D6").Value = "Lelouch" + [char]10 + "Code Geass"

# --- Row 7: Hat image link ---
$ws.Range("G7").Value = "https://cdn.discordapp.com/attachments/699111007649398865/1048523933063843931/Dunk_Sweatling_Lelouch_HatHair_V2b_210x210.png"

# --- Row 9: Neck/outfit image link ---
$ws.Range("G9").Value = "https://cdn.discordapp.com/attachments/699111007649398865/1048523933399404614/Dunk_Sweatling_Lelouch_NeckClothes_V2b_210x210.png"

# --- Row 12: creator lookup info ---
$ws.Range("C12").Value = "omnipotent_0"
$ws.Range("D12").Value = 42256416
$ws.Range("E12").Value = "257 days"

# --- Insert a new row 15 (old row 15 "21 / 1081" shifts down to row 16) ---
$ws.Rows.Item(15).Insert()

# Row 14 used to hold the "shop item rows / citb user(s) / citb comment"
# footer labels; that footer now lives on the newly inserted row 15, so
# row 14 becomes a normal (empty in L:N) redeemer row instead.
$ws.Range("A14").Value = 19
$ws.Range("C14").Value = "Omnipotent_0"
$ws.Range("L14").Value = "'"
$ws.Range("M14").Value = "'"
$ws.Range("N14").Value = "'"

# New row 15: blank data columns, with the footer labels moved here
$ws.Range("A15").Value = 20
$ws.Range("B15").Value = "'"
$ws.Range("C15").Value = "'"
$ws.Range("D15").Value = "'"
$ws.Range("E15").Value = "'"
$ws.Range("F15").Value = "'"
$ws.Range("G15").Value = "'"
$ws.Range("H15").Value = "'"
$ws.Range("I15").Value = "'"
$ws.Range("J15").Value = "'"
$ws.Range("K15").Value = "'"
$ws.Range("L15").Value = "shop item rows"
$ws.Range("M15").Value = "citb user(s)"
$ws.Range("N15").Value = "citb comment"

# Row 16 (previously row 15): shop-row counter now 822 instead of 1081
$ws.Range("L16").Value = 822
